$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, matching style of neighboring header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timeTaken = @(
  "2021-10-05 13:40:17.814510",
  "2021-10-05 13:40:17.814522",
  "2021-10-05 13:40:17.814526",
  "2021-10-05 13:40:17.814528",
  "2021-10-05 13:40:17.814531",
  "2021-10-05 13:40:17.814534",
  "2021-10-05 13:40:17.814537",
  "2021-10-05 13:40:17.814539",
  "2021-10-05 13:40:17.814542",
  "2021-10-05 13:40:17.814545",
  "2021-10-05 13:40:17.814548",
  "2021-10-05 13:40:17.814550",
  "2021-10-05 13:40:17.814553",
  "2021-10-05 13:40:17.814556",
  "2021-10-05 13:40:17.814558",
  "2021-10-05 13:40:17.814561",
  "2021-10-05 13:40:17.814563",
  "2021-10-05 13:40:17.814566",
  "2021-10-05 13:40:17.814569",
  "2021-10-05 13:40:17.814571",
  "2021-10-05 13:40:17.814574",
  "2021-10-05 13:40:17.814576",
  "2021-10-05 13:40:17.814579",
  "2021-10-05 13:40:17.814582",
  "2021-10-05 13:40:17.814584",
  "2021-10-05 13:40:17.814587",
  "2021-10-05 13:40:17.814590",
  "2021-10-05 13:40:17.814593",
  "2021-10-05 13:40:17.814595",
  "2021-10-05 13:40:17.814598",
  "2021-10-05 13:40:17.814600",
  "2021-10-05 13:40:17.814603",
  "2021-10-05 13:40:17.814606",
  "2021-10-05 13:40:17.814609",
  "2021-10-05 13:40:17.814611",
  "2021-10-05 13:40:17.814614",
  "2021-10-05 13:40:17.814616",
  "2021-10-05 13:40:17.814619",
  "2021-10-05 13:40:17.814621",
  "2021-10-05 13:40:17.814624",
  "2021-10-05 13:40:17.814627",
  "2021-10-05 13:40:17.814630",
  "2021-10-05 13:40:17.814632",
  "2021-10-05 13:40:17.814635",
  "2021-10-05 13:40:17.814637",
  "2021-10-05 13:40:17.814640",
  "2021-10-05 13:40:17.814643",
  "2021-10-05 13:40:17.814645",
  "2021-10-05 13:40:17.814648",
  "2021-10-05 13:40:17.814650",
  "2021-10-05 13:40:17.814653",
  "2021-10-05 13:40:17.814655",
  "2021-10-05 13:40:17.814658",
  "2021-10-05 13:40:17.814661",
  "2021-10-05 13:40:17.814664",
  "2021-10-05 13:40:17.814666",
  "2021-10-05 13:40:17.814669",
  "2021-10-05 13:40:17.814671",
  "2021-10-05 13:40:17.814674",
  "2021-10-05 13:40:17.814677",
  "2021-10-05 13:40:17.814679",
  "2021-10-05 13:40:17.814682",
  "2021-10-05 13:40:17.814684",
  "2021-10-05 13:40:17.814687",
  "2021-10-05 13:40:17.814691",
  "2021-10-05 13:40:17.814694",
  "2021-10-05 13:40:17.814697",
  "2021-10-05 13:40:17.814700",
  "2021-10-05 13:40:17.814702",
  "2021-10-05 13:40:17.814705",
  "2021-10-05 13:40:17.814708",
  "2021-10-05 13:40:17.814710",
  "2021-10-05 13:40:17.814713",
  "2021-10-05 13:40:17.814716",
  "2021-10-05 13:40:17.814718",
  "2021-10-05 13:40:17.814721",
  "2021-10-05 13:40:17.814725",
  "2021-10-05 13:40:17.814728",
  "2021-10-05 13:40:17.814731",
  "2021-10-05 13:40:17.814734",
  "2021-10-05 13:40:17.814736",
  "2021-10-05 13:40:17.814739",
  "2021-10-05 13:40:17.814742",
  "2021-10-05 13:40:17.814744",
  "2021-10-05 13:40:17.814747",
  "2021-10-05 13:40:17.814750",
  "2021-10-05 13:40:17.814752",
  "2021-10-05 13:40:17.814755",
  "2021-10-05 13:40:17.814758",
  "2021-10-05 13:40:17.814760",
  "2021-10-05 13:40:17.814763",
  "2021-10-05 13:40:17.814766",
  "2021-10-05 13:40:17.814770",
  "2021-10-05 13:40:17.814773",
  "2021-10-05 13:40:17.814775",
  "2021-10-05 13:40:17.814778",
  "2021-10-05 13:40:17.814781",
  "2021-10-05 13:40:17.814783",
  "2021-10-05 13:40:17.814786",
  "2021-10-05 13:40:17.814789",
  "2021-10-05 13:40:17.814791",
  "2021-10-05 13:40:17.814794",
  "2021-10-05 13:40:17.814797",
  "2021-10-05 13:40:17.814799",
  "2021-10-05 13:40:17.814802",
  "2021-10-05 13:40:17.814805",
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Output "done"